# Fruta / hortaliza, semanal
# Insert a new weekly data row at row 468 (pushing existing rows 468:570
# down to 469:571) and populate it with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 468..570 down by one row, creating a blank row 468.
$ws.Rows("468:468").Insert()

# Populate the newly inserted row with the new record.
$ws.Range("A468").Value = 3
$ws.Range("B468").Value = "Femacal de La Calera"
$ws.Range("C468").Value = "Coquimbo"
$ws.Range("D468").Value = 44889
$ws.Range("E468").Value = 5
$ws.Range("F468").Value = 100112021
$ws.Range("G468").Value = "Ají"
$ws.Range("H468").Value = "Inferno"
$ws.Range("I468").Value = "Primera"
$ws.Range("J468").Value = 76
$ws.Range("K468").Value = 20000
$ws.Range("L468").Value = 21000
$ws.Range("M468").Value = 20500
$ws.Range("N468").Value = "`$/caja 15 kilos"
$ws.Range("O468").Value = "Limache"
$ws.Range("P468").Value = 1367
$ws.Range("Q468").Value = 15
$ws.Range("R468").Value = "Hortaliza"
